$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("C3") "1014"
Set-TextValue $ws.Range("D3") "3237394.33"
Set-TextValue $ws.Range("C9") "59"
Set-TextValue $ws.Range("D9") "155597.64"
Set-TextValue $ws.Range("C10") "364"
Set-TextValue $ws.Range("D10") "1303761.71"
Set-TextValue $ws.Range("C11") "149"
Set-TextValue $ws.Range("D11") "647101.77"
Set-TextValue $ws.Range("C12") "36"
Set-TextValue $ws.Range("D12") "183783.00"
Set-TextValue $ws.Range("C16") "428"
Set-TextValue $ws.Range("D16") "1357961.23"
Set-TextValue $ws.Range("C17") "155"
Set-TextValue $ws.Range("D17") "671702.10"
Set-TextValue $ws.Range("C18") "48"
Set-TextValue $ws.Range("D18") "251045.00"
Set-TextValue $ws.Range("C22") "327"
Set-TextValue $ws.Range("D22") "963679.20"
Set-TextValue $ws.Range("C34") "574"
Set-TextValue $ws.Range("D34") "1891313.66"
Set-TextValue $ws.Range("C35") "229"
Set-TextValue $ws.Range("D35") "1157788.11"
Set-TextValue $ws.Range("C51") "104"
Set-TextValue $ws.Range("D51") "297768.17"
Set-TextValue $ws.Range("C52") "597"
Set-TextValue $ws.Range("D52") "2107095.21"
Set-TextValue $ws.Range("C53") "265"
Set-TextValue $ws.Range("D53") "1172878.76"
Set-TextValue $ws.Range("C54") "88"
Set-TextValue $ws.Range("D54") "511378.23"
Set-TextValue $ws.Range("C56") "24"
Set-TextValue $ws.Range("D56") "78220.65"
Set-TextValue $ws.Range("C67") "12"
Set-TextValue $ws.Range("D67") "64027.00"
